# Apply edits described by the diff: title rewording, "What we like" /
# "What we don't like" bullet rewording, and the closing bold/italic
# summary rewording.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Title (appears twice: the Heading1 at the top, and the bold run near
# the bottom) - ReplaceAll via Content.Find handles both occurrences.
Replace-Text "Play Crystal Forest HD for Free - WMS Online Slot Review" `
             "Play Crystal Forest HD Free - Exciting Slot Game Review"

# "What we like" bullets
Replace-Text "Enchanting, detailed graphics" `
             "Traditional structure with 5 reels and 25 paylines"

Replace-Text "Free spins feature" `
             "Free spins and cascading reels add excitement to gameplay"

Replace-Text "Cascading reels" `
             "Enchanting and detailed graphics create a magical atmosphere"

Replace-Text "Playable on desktop and mobile" `
             "Available on both desktop and mobile devices"

# "What we don't like" bullets
Replace-Text "No scatter symbol" `
             "Limited variety of bonus features"

Replace-Text "Limited bonus features" `
             "No scatter symbol in the game"

# Closing italic summary paragraph
Replace-Text "Get familiar with the tricks and features of WMS slot game Crystal Forest HD. Play for free and enjoy enchanting graphics and cascading reels." `
             "Read our review of Crystal Forest HD, a slot game with free spins and enchanting graphics. Play for free!"
